$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = "Gross Profit"
$ws.Range("C7").Value = "Gross Profit"
$ws.Range("C8").Value = "Gross Profit"
$ws.Range("C9").Value = "Gross Profit"
$ws.Range("E9").Value = "Operational Income"

$ws.Range("E5").Select()
